$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (so old F becomes G)
$ws.Columns("F").Insert()

# Set header for new column F
$ws.Range("F1").Value = "anotaciones"

# Set value for new column F row 2
$ws.Range("F2").Value = "soyyo"

# Update the date value in the (now) G2 cell, keeping it as text (not auto-converted to a date)
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2023-09-07"
